# "save the PAM data!" - add a `species` column (M) to the PAM fluorometry
# data sheet, labeling every data row with the anemone cross species name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$speciesHeader = "species"
$speciesValue  = "A. elegantissima x B. muscatinei "

# Header for the new column.
$ws.Range("M1").Value = $speciesHeader

# Fill the species value down every data row (row 2 is a sub-header row and
# stays untouched, matching the source data).
for ($r = 3; $r -le 62; $r++) {
    $ws.Cells.Item($r, 13).Value = $speciesValue
}

# Reflect the scrolled/selected state that was saved along with the new data.
$win = $excel.ActiveWindow
$win.ScrollRow = 36
$win.ScrollColumn = 1
$ws.Range("M3:M62").Select()
